$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.947.08"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.636.49"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.80"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  -0.80%  "
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.32"
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("D13").Value = "1.637.63"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.95"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("E16").Value = "  -1.57%  "
$ws.Range("D17").Value = "25.844.26"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "192.78"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.12"
$ws.Range("E22").Value = "  -2.57%  "
$ws.Range("E23").Value = "  +2.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.68"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.76"
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.42"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  -2.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.24"
$ws.Range("E32").Value = "  -1.59%  "
$ws.Range("E33").Value = "  -1.38%  "
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.893"
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("D36").Value = "1.121.03"
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.46"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").Value = "  -2.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0156"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.794"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.09"
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.31"
$ws.Range("E42").Value = "  -3.29%  "
$ws.Range("D43").Value = "0.0₆0114"
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0520"
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.47"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.70"
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.414"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0937"
$ws.Range("E50").Value = "  -3.17%  "
$ws.Range("E51").Value = "  -1.94%  "
